$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SKU code update (A3: 74533 -> 70286)
$ws.Range("A3").Value = 70286

# Price column D3:D5 switch from numeric currency values to literal Thai-baht
# priced text strings (stored as shared strings once saved).
$ws.Range("D3").Value = "฿279.00"
$ws.Range("D4").Value = "฿129.00"
$ws.Range("D5").Value = "฿79.00"

# D7 price text corrected to include the trailing zero.
$ws.Range("D7").Value = "฿74.00"

# Move the active selection to D14.
$ws.Activate()
$ws.Range("D14").Select()

$wb.Save()
